$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2:B7 from "No" to "Yes"
$ws.Range("B2:B7").Value = "Yes"

# Set the active selection to C4
$ws.Range("C4").Select()
